$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the test case previously labeled TC_011 (row 12) to TC_000
$ws.Range("A12").Value = "TC_000"

# Insert new row data for TC_011 (globalContact.feature - launch browser)
$ws.Range("A13").Value = "TC_011"
$ws.Range("B13").Value = "globalContact.feature"
$ws.Range("C13").Value = "Launch Browser and go to application"
$ws.Range("D13").Value = "No"
$ws.Range("E13").Value = "No"

# Insert new row data for TC_012 (globalContact.feature - create individual contact)
$ws.Range("A14").Value = "TC_012"
$ws.Range("B14").Value = "globalContact.feature"
$ws.Range("C14").Value = "Verify user enters first and last name then clicks Create Individual Contact and lands on the Individual Contact page with pre-filled fields"
$ws.Range("D14").Value = "No"
$ws.Range("E14").Value = "No"

# Update SmokeTest column values that were "Yes" to "No"
$ws.Range("D2").Value = "No"
$ws.Range("D12").Value = "No"

# Update the active selection cell from C15 to D15
$ws.Range("D15").Select()

# Extend the SmokeTest/RegressionTest list validation to cover the new rows (D13:D14)
$ws.Range("D1:D14").Validation.Modify(3, 1, 1, '"Yes,No"')
